$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: B/C change from the old "Objetivos" paragraph to the
#     "519033 - Carlos Yujiro Shigue" docente string (cells already exist,
#     so a plain value write keeps the existing style). ---
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

# --- Row 13: becomes "Programa resumido:" / "01/01/2012" / "01/01/2012",
#     with a new A13 cell (style copied from another column-A cell) and
#     row height 60. The "01/01/2012" text would auto-parse as a date via
#     a plain .Value assignment, so it is written as a formula returning
#     text and then converted back to a plain value (paste-values) to
#     land as a shared string with no date number-format. ---
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").Value = "Programa resumido:"

$ws.Range("B13").Formula = '="01/01/2012"'
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4163) | Out-Null

$ws.Range("C13").Formula = '="01/01/2012"'
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4163) | Out-Null

$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: becomes "Short syllabus:" only, row height 60. B14/C14
#     (previously holding the Katia docente string) must no longer exist. ---
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Clear() | Out-Null
$ws.Range("C14").Clear() | Out-Null
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15: becomes "Programa:" / "519033 - Carlos Yujiro Shigue" (x2),
#     row height 120. ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C15").Value = "519033 - Carlos Yujiro Shigue"
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16: becomes "Syllabus:" only, row height 120. ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17: becomes "Avaliação:" only, default row height (was 120).
#     B17/C17 (old "Programa:" paragraph) must no longer exist. ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Clear() | Out-Null
$ws.Range("C17").Clear() | Out-Null
$ws.Rows.Item(17).AutoFit() | Out-Null

# --- Row 18: becomes "Método:" / "5817692 - Katia Cristiane Gandolpho
#     Candioto" (x2). B18/C18 are new cells, so copy format from B15/C15
#     (style 2 / style 3) before writing the value. Row height 60. ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B15").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("B18").Value = "5817692 - Katia Cristiane Gandolpho Candioto"
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Value = "5817692 - Katia Cristiane Gandolpho Candioto"
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19: becomes "Critério:" / the "Aulas expositivas..." paragraph
#     (x2). B19/C19 are new cells too. Row height 60. ---
$ws.Range("A19").Value = "Critério:"
$ws.Range("B15").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null
$ws.Range("B19").Value = "Aulas expositivas, reuniões com professor orientador, desenvolvimento de projeto de pesquisa e elaboração de monografia."
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$ws.Range("C19").Value = "Aulas expositivas, reuniões com professor orientador, desenvolvimento de projeto de pesquisa e elaboração de monografia."
$ws.Rows.Item(19).RowHeight = 60

# --- Row 20: becomes "Norma de recuperação:" / "Nota de avaliação da
#     monografia." (x2), row height 60 (cells already exist). ---
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Nota de avaliação da monografia."
$ws.Range("C20").Value = "Nota de avaliação da monografia."
$ws.Rows.Item(20).RowHeight = 60

# --- Row 21: becomes "Bibliografia:" / "A critério da Comissão de Curso
#     poderá ser oferecida recuperação." (x2), row height 120. ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "A critério da Comissão de Curso poderá ser oferecida recuperação."
$ws.Range("C21").Value = "A critério da Comissão de Curso poderá ser oferecida recuperação."
$ws.Rows.Item(21).RowHeight = 120

# --- Row 22: becomes "Requisitos:" only, default row height (was 60).
#     B22/C22 (the old long bibliography paragraph) must no longer exist. ---
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Clear() | Out-Null
$ws.Range("C22").Clear() | Out-Null
$ws.Rows.Item(22).AutoFit() | Out-Null

# --- Row 23: becomes the LOB1008 requirement line (x2) only (no A23),
#     row height 30. A23 (old "Bibliografia:" label) must no longer exist. ---
$ws.Range("A23").Clear() | Out-Null
$ws.Range("B23").Value = "LOB1008 -  Ciência, Tecnologia e Sociedade  (Requisito)`n"
$ws.Range("C23").Value = "LOB1008 -  Ciência, Tecnologia e Sociedade  (Requisito)`n"
$ws.Rows.Item(23).RowHeight = 30

# --- Row 24: becomes the LOB1045 requirement line (x2) only (no A24),
#     row height 30. B24/C24 are new cells here. ---
$ws.Range("A24").Clear() | Out-Null
$ws.Range("B23").Copy() | Out-Null
$ws.Range("B24").PasteSpecial(-4122) | Out-Null
$ws.Range("B24").Value = "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)`n"
$ws.Range("C23").Copy() | Out-Null
$ws.Range("C24").PasteSpecial(-4122) | Out-Null
$ws.Range("C24").Value = "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)`n"
$ws.Rows.Item(24).RowHeight = 30

# --- Rows 25 and 26 (old Bibliografia requirement lines) no longer
#     exist; remove them so the sheet shrinks to A1:C24. ---
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(25).Delete()
